# Update cryptocurrency price / 1h-volume figures in the active worksheet
# to reflect the latest scrape (GitHub Actions cron update).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell's value while forcing it to remain plain text (so that
# numeric-looking strings like "1.00" or "57.279.97" are not reinterpreted
# by Excel as numbers and keep their exact displayed formatting), then drop
# the temporary "Text" number format again so no stray cell styling is left
# behind on the cell.
function Set-TextValue {
    param($Address, $Text)
    $cell = $ws.Range($Address)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

Set-TextValue 'D2' '57.279.97'
$ws.Range('E2').Value = '  -1.02%  '
Set-TextValue 'D3' '3.104.09'
$ws.Range('E3').Value = '  -0.07%  '
Set-TextValue 'D4' '1.00'
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '525.95'
$ws.Range('E5').Value = '  +0.07%  '
Set-TextValue 'D6' '137.49'
$ws.Range('E6').Value = '  -3.46%  '
Set-TextValue 'D7' '1.00'
$ws.Range('E7').Value = '  +0.02%  '
Set-TextValue 'D8' '3.101.55'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('E9').Value = '  +2.34%  '
Set-TextValue 'D10' '7.34'
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('E11').Value = '  -1.04%  '
Set-TextValue 'D12' '0.401'
$ws.Range('E12').Value = '  +2.26%  '
Set-TextValue 'D13' '3.633.42'
$ws.Range('E13').Value = '  -0.14%  '
Set-TextValue 'D14' '0.135'
$ws.Range('E14').Value = '  +1.61%  '
Set-TextValue 'D15' '25.57'
$ws.Range('E15').Value = '  -0.31%  '
$ws.Range('E16').Value = '  -1.28%  '
Set-TextValue 'D17' '57.379.13'
$ws.Range('E17').Value = '  -1.00%  '
Set-TextValue 'D18' '3.103.52'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('E19').Value = '  -2.82%  '
Set-TextValue 'D20' '12.54'
$ws.Range('E20').Value = '  -1.83%  '
$ws.Range('E21').Value = '  -0.86%  '
Set-TextValue 'D22' '350.05'
$ws.Range('E22').Value = '  +2.50%  '
$ws.Range('E23').Value = '  +0.06%  '
Set-TextValue 'D24' '68.10'
$ws.Range('E24').Value = '  +1.28%  '
Set-TextValue 'D25' '0.502'
$ws.Range('E25').Value = '  -2.20%  '
$ws.Range('E26').Value = '  -1.26%  '
$ws.Range('E27').Value = '  -0.19%  '
Set-TextValue 'D28' '0.0₃0892'
$ws.Range('E28').Value = '  -2.89%  '
Set-TextValue 'D29' '0.999'
$ws.Range('E29').Value = '  +0.03%  '
Set-TextValue 'D30' '7.35'
$ws.Range('E30').Value = '  +1.03%  '
$ws.Range('E31').Value = '  +0.33%  '
Set-TextValue 'D32' '5.98'
$ws.Range('E32').Value = '  -7.47%  '
$ws.Range('E33').Value = '  -0.66%  '
$ws.Range('E34').Value = '  +7.87%  '
$ws.Range('E35').Value = '  -3.94%  '
Set-TextValue 'D36' '159.41'
$ws.Range('E36').Value = '  +0.85%  '
$ws.Range('E37').Value = '  -1.66%  '
Set-TextValue 'D38' '26.36'
$ws.Range('E38').Value = '  -0.24%  '
Set-TextValue 'D39' '1.25'
$ws.Range('E39').Value = '  -1.05%  '
Set-TextValue 'D40' '0.0658'
$ws.Range('E40').Value = '  -1.70%  '
$ws.Range('E41').Value = '  +0.92%  '
$ws.Range('E42').Value = '  +1.11%  '
Set-TextValue 'D43' '0.695'
$ws.Range('E43').Value = '  +1.62%  '
Set-TextValue 'D44' '2.407.39'
$ws.Range('E44').Value = '  +5.22%  '
Set-TextValue 'D45' '36.69'
$ws.Range('E45').Value = '  -0.44%  '
Set-TextValue 'D46' '1.00'
$ws.Range('E46').Value = '  -0.05%  '
Set-TextValue 'D47' '3.140.87'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  +0.67%  '
Set-TextValue 'D49' '0.970'
$ws.Range('E49').Value = '  -3.05%  '
$ws.Range('E50').Value = '  -1.70%  '
Set-TextValue 'D51' '0.765'
$ws.Range('E51').Value = '  +2.47%  '
